$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated stock/price values (row r => B/D column updates)
$ws.Range("B2").Value = 71429.02
$ws.Range("D2").Value = 1197.5
$ws.Range("B5").Value = 74800
$ws.Range("B6").Value = 317000
$ws.Range("B8").Value = 79900
$ws.Range("B9").Value = 123000
$ws.Range("B10").Value = 78300
$ws.Range("B11").Value = 8165
$ws.Range("B12").Value = 18600
$ws.Range("B13").Value = 11435
$ws.Range("B14").Value = 28630
$ws.Range("B15").Value = 23385
$ws.Range("B17").Value = 5965
$ws.Range("B19").Value = 13445
$ws.Range("B20").Value = 15270
$ws.Range("B21").Value = 455500
$ws.Range("B22").Value = 54700
$ws.Range("B32").Value = 53379000
$ws.Range("B33").Value = 3791000
$ws.Range("B34").Value = 2.1310345312666299

# Scroll/selection state update to match the saved view
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A26").Select()
